# Sample Project / Main.xlsx — "1" save edit
#
# The rule table's last row (row 11, the "R40" rule) has its Rule-name
# cell (B11) changed from the text "R40" to the text "1". The cell keeps
# its existing style/format (s="23") and keeps being stored as a text
# (shared-string) value rather than turning into a number.
#
# Setting .Value = "1" directly would get auto-coerced to a numeric
# literal by Excel's type inference (since "1" parses as a number), which
# would also drop the t="s" string-cell marker. To force it to remain
# text without touching the cell's number format/style, we write it as a
# text formula and then convert that formula to its literal value via
# copy / paste-special-values (xlPasteValues = -4163), exactly like using
# Excel's own "Convert to values" command.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)
